# Auto-update gym prices
$wb = $excel.ActiveWorkbook

# "4x4 Squat Racks" sheet: update price in C2
$ws1 = $wb.Worksheets.Item("4x4 Squat Racks")
$ws1.Range("C2").Value = "$2,131.00"

# "Squat Stands" sheet: update price in C2, clear price in C3
$ws2 = $wb.Worksheets.Item("Squat Stands")
$ws2.Range("C2").Value = "$1,541.00"
$ws2.Range("C3").ClearContents()
